$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns stay text (workbook stores them as inline strings,
# e.g. "43.783.96" / "1.00" / "  -0.22%  " -- without this Excel would coerce
# numeric-looking values to numbers and drop formatting such as trailing zeros).
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '43.783.96'
$ws.Range("E2").Value = '  -0.22%  '
# Row 3
$ws.Range("D3").Value = '2.291.03'
$ws.Range("E3").Value = '  -0.22%  '
# Row 4
$ws.Range("E4").Value = '  -0.15%  '
# Row 5
$ws.Range("D5").Value = '123.16'
$ws.Range("E5").Value = '  +8.14%  '
# Row 6
$ws.Range("D6").Value = '267.61'
$ws.Range("E6").Value = '  -0.95%  '
# Row 7
$ws.Range("E7").Value = '  +2.00%  '
# Row 8
$ws.Range("E8").Value = '  +0.06%  '
# Row 9
$ws.Range("D9").Value = '0.621'
$ws.Range("E9").Value = '  +0.19%  '
# Row 10
$ws.Range("D10").Value = '48.92'
$ws.Range("E10").Value = '  +1.60%  '
# Row 11
$ws.Range("E11").Value = '  -0.20%  '
# Row 12
$ws.Range("D12").Value = '9.18'
$ws.Range("E12").Value = '  +1.30%  '
# Row 13
$ws.Range("E13").Value = '  +0.81%  '
# Row 14
$ws.Range("D14").Value = '15.62'
$ws.Range("E14").Value = '  -1.78%  '
# Row 15
$ws.Range("D15").Value = '0.895'
$ws.Range("E15").Value = '  +4.01%  '
# Row 16
$ws.Range("D16").Value = '2.635.73'
$ws.Range("E16").Value = '  -0.15%  '
# Row 17
$ws.Range("D17").Value = '2.292.13'
$ws.Range("E17").Value = '  -0.13%  '
# Row 18
$ws.Range("D18").Value = '43.758.19'
$ws.Range("E18").Value = '  +0.01%  '
# Row 19
$ws.Range("E19").Value = '  +0.00%  '
# Row 20
$ws.Range("D20").Value = '7.06'
$ws.Range("E20").Value = '  +2.59%  '
# Row 21
$ws.Range("D21").Value = '72.63'
$ws.Range("E21").Value = '  +0.61%  '
# Row 22
$ws.Range("D22").Value = '2.45'
$ws.Range("E22").Value = '  +0.75%  '
# Row 23
$ws.Range("D23").Value = '236.22'
$ws.Range("E23").Value = '  +1.38%  '
# Row 24
$ws.Range("D24").Value = '9.74'
$ws.Range("E24").Value = '  +0.59%  '
# Row 25
$ws.Range("D25").Value = '2.89'
$ws.Range("E25").Value = '  -4.42%  '
# Row 26
$ws.Range("E26").Value = '  +1.65%  '
# Row 27
$ws.Range("D27").Value = '11.82'
$ws.Range("E27").Value = '  +1.47%  '
# Row 28
$ws.Range("D28").Value = '42.78'
$ws.Range("E28").Value = '  +2.16%  '
# Row 29
$ws.Range("E29").Value = '  +0.38%  '
# Row 30
$ws.Range("E30").Value = '  -0.52%  '
# Row 31
$ws.Range("D31").Value = '173.44'
$ws.Range("E31").Value = '  -1.24%  '
# Row 32
$ws.Range("D32").Value = '21.74'
$ws.Range("E32").Value = '  +0.68%  '
# Row 33
$ws.Range("D33").Value = '0.0914'
$ws.Range("E33").Value = '  -1.27%  '
# Row 34
$ws.Range("E34").Value = '  +1.49%  '
# Row 35
$ws.Range("E35").Value = '  +2.15%  '
# Row 36
$ws.Range("D36").Value = '0.0381'
$ws.Range("E36").Value = '  +4.52%  '
# Row 37
$ws.Range("D37").Value = '4.74'
$ws.Range("E37").Value = '  +1.40%  '
# Row 38
$ws.Range("D38").Value = '4.05'
$ws.Range("E38").Value = '  +5.21%  '
# Row 39
$ws.Range("D39").Value = '0.108'
$ws.Range("E39").Value = '  +0.71%  '
# Row 40
$ws.Range("D40").Value = '2.57'
$ws.Range("E40").Value = '  +7.78%  '
# Row 41
$ws.Range("D41").Value = '14.41'
$ws.Range("E41").Value = '  +3.53%  '
# Row 42
$ws.Range("D42").Value = '75.34'
$ws.Range("E42").Value = '  +1.89%  '
# Row 43
$ws.Range("D43").Value = '0.240'
$ws.Range("E43").Value = '  -0.83%  '
# Row 44
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  -0.06%  '
# Row 45
$ws.Range("B45").Value = 'THORChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D45").Value = '5.97'
$ws.Range("E45").Value = '  -5.45%  '
# Row 46
$ws.Range("E46").Value = '  -0.83%  '
# Row 47
$ws.Range("E47").Value = '  +3.18%  '
# Row 48
$ws.Range("B48").Value = 'ordi'
$ws.Range("C48").Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range("D48").Value = '73.74'
$ws.Range("E48").Value = '  +37.15%  '
# Row 49
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '8.59'
$ws.Range("E49").Value = '  -2.78%  '
# Row 50
$ws.Range("E50").Value = '  +0.81%  '
# Row 51
$ws.Range("D51").Value = '102.10'
$ws.Range("E51").Value = '  +0.57%  '
